$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value (as scraped from the updated
# coinranking.com feed). Values must stay plain text (e.g. "305.91",
# "-0.54%") exactly as the source XML stores them as inline strings, so we
# force text formatting before assigning and then drop back to the default
# "Normal" style so no stray formatting is introduced.
$updates = @{
    "D2" = "305.91"
    "E2" = "-0.54%"
    "D3" = "38.90"
    "E3" = "7.35%"
    "D4" = "5.113"
    "E4" = "1.00%"
    "D5" = "0.08074"
    "E5" = "-0.37%"
    "D6" = "1.928"
    "E6" = "-2.70%"
    "D7" = "4.198"
    "E7" = "0.70%"
    "D8" = "7.997"
    "E8" = "1.58%"
    "D9" = "0.9281"
    "E9" = "0.06%"
    "D10" = "0.1479"
    "E10" = "1.67%"
    "D11" = "0.1912"
    "E11" = "-1.17%"
    "D12" = "0.09154"
    "E12" = "0.59%"
    "D13" = "0.03517"
    "E13" = "2.24%"
    "D14" = "0.09773"
    "E14" = "-1.11%"
    "D15" = "0.001399"
    "E15" = "-0.56%"
    "D16" = "0.006053"
    "E16" = "-5.35%"
    "D17" = "3.780"
    "E17" = "-1.42%"
    "D18" = "3.411"
    "E18" = "-0.57%"
    "D19" = "0.3431"
    "E19" = "-0.70%"
    "D20" = "0.1321"
    "E20" = "-0.01%"
    "D21" = "4.684"
    "E21" = "-2.74%"
    "D22" = "0.2421"
    "E22" = "3.18%"
    "D23" = "0.04385"
    "E23" = "-0.14%"
    "D24" = "0.001238"
    "E24" = "0.12%"
    "D25" = "0.004264"
    "E25" = "2.16%"
    "D26" = "0.0001303"
    "E26" = "0.00%"
    "D39" = "0.02028"
    "E39" = "-0.52%"
    "D40" = "0.05044"
    "E40" = "-1.54%"
    "D41" = "0.007516"
    "E41" = "0.65%"
    "D42" = "0.009720"
    "E42" = "-3.40%"
    "E43" = "-1.95%"
    "D44" = "0.002104"
    "E44" = "-0.95%"
    "D45" = "0.009903"
    "E45" = "0.39%"
    "D46" = "0.00006198"
    "E46" = "-1.85%"
    "E47" = "-0.01%"
    "D48" = "0.002873"
    "D49" = "0.001806"
    "E49" = "12.61%"
    "D50" = "0.00002104"
    "E50" = "-0.01%"
    "D51" = "0.0002004"
    "E51" = "-0.01%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
